$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D (Price) to Text format so numeric-looking
# strings (e.g. "0.6374") are stored as text, matching the source data,
# then restore the original (default) style once all values are set.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.884.93"
$ws.Range("E2").Value = "  +2.73%  "
$ws.Range("D3").Value = "1.858.82"
$ws.Range("E3").Value = "  +2.05%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "246.37"
$ws.Range("E5").Value = "  +2.12%  "
$ws.Range("D6").Value = "0.6374"
$ws.Range("E6").Value = "  +3.61%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "0.2999"
$ws.Range("E8").Value = "  +3.78%  "
$ws.Range("D9").Value = "0.07470"
$ws.Range("E9").Value = "  +1.97%  "
$ws.Range("D10").Value = "24.52"
$ws.Range("E10").Value = "  +7.05%  "
$ws.Range("D11").Value = "0.07677"
$ws.Range("E11").Value = "  +0.26%  "
$ws.Range("D12").Value = "1.882.77"
$ws.Range("E12").Value = "  +3.35%  "
$ws.Range("D13").Value = "5.049"
$ws.Range("E13").Value = "  +2.14%  "
$ws.Range("D14").Value = "0.6913"
$ws.Range("E14").Value = "  +4.92%  "
$ws.Range("D15").Value = "84.31"
$ws.Range("E15").Value = "  +3.21%  "
$ws.Range("D16").Value = "0.000009336"
$ws.Range("E16").Value = "  +4.05%  "
$ws.Range("D17").Value = "6.073"
$ws.Range("E17").Value = "  +4.31%  "
$ws.Range("D18").Value = "29.843.91"
$ws.Range("E18").Value = "  +2.71%  "
$ws.Range("D19").Value = "2.117.03"
$ws.Range("E19").Value = "  +2.58%  "
$ws.Range("D20").Value = "238.54"
$ws.Range("E20").Value = "  +0.29%  "
$ws.Range("D21").Value = "12.66"
$ws.Range("E21").Value = "  +1.91%  "
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").Value = "7.360"
$ws.Range("E23").Value = "  +3.59%  "
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").Value = "159.13"
$ws.Range("E25").Value = "  +1.19%  "
$ws.Range("D26").Value = "0.1416"
$ws.Range("E26").Value = "  +0.84%  "
$ws.Range("D27").Value = "8.578"
$ws.Range("E27").Value = "  +2.07%  "
$ws.Range("D28").Value = "17.97"
$ws.Range("E28").Value = "  +2.13%  "
$ws.Range("E29").Value = "  +1.61%  "
$ws.Range("D30").Value = "0.06058"
$ws.Range("E30").Value = "  +9.16%  "
$ws.Range("D31").Value = "1.275"
$ws.Range("E31").Value = "  +5.74%  "
$ws.Range("D32").Value = "4.131"
$ws.Range("E32").Value = "  +0.98%  "
$ws.Range("E33").Value = "  +1.29%  "
$ws.Range("D34").Value = "1.891"
$ws.Range("E34").Value = "  +4.24%  "
$ws.Range("E35").Value = "  +3.23%  "
$ws.Range("D36").Value = "0.7280"
$ws.Range("E36").Value = "  -0.76%  "
$ws.Range("D37").Value = "2.611"
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("D38").Value = "2.857"
$ws.Range("E38").Value = "  +1.06%  "
$ws.Range("E39").Value = "  +2.47%  "
$ws.Range("D40").Value = "1.223.53"
$ws.Range("E40").Value = "  +1.27%  "
$ws.Range("D41").Value = "0.9351"
$ws.Range("E41").Value = "  +4.66%  "
$ws.Range("D42").Value = "6.285"
$ws.Range("E42").Value = "  -1.03%  "
$ws.Range("D43").Value = "1.002"
$ws.Range("E43").Value = "  +0.33%  "
$ws.Range("D44").Value = "2.023.51"
$ws.Range("E44").Value = "  +2.71%  "
$ws.Range("D45").Value = "102.28"
$ws.Range("E45").Value = "  +1.42%  "
$ws.Range("D46").Value = "66.25"
$ws.Range("E46").Value = "  +2.64%  "
$ws.Range("D47").Value = "0.00000000122"
$ws.Range("E47").Value = "  +4.55%  "
$ws.Range("D48").Value = "0.5090"
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("D49").Value = "9.275"
$ws.Range("E49").Value = "  +2.63%  "
$ws.Range("E50").Value = "  +2.35%  "
$ws.Range("D51").Value = "0.1141"
$ws.Range("E51").Value = "  +2.94%  "

$ws.Range("D2:D51").Style = "Normal"
